# Update "Mata Kuliah" template: add a new "ID_SMS" column right after the
# first ("No") column on the data sheet, then restore the expected view
# state (zoom, selection, active sheet).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Mata Kuliah"
$ws2 = $wb.Worksheets.Item(2)   # "Petunjuk Pengisian"

# Insert a new column before column B; this shifts the existing C:S data
# right to D:T and carries the column-A formatting (style) into the new
# column B, matching how the header/data row styles line up afterwards.
$ws1.Columns.Item(2).Insert()

# New header cell for the inserted column.
$ws1.Range("B1").Value = "ID_SMS"

# "Petunjuk Pengisian" keeps its own last selection / scroll position.
$ws2.Activate()
$ws2.Range("F7").Select()

# Restore view state: "Mata Kuliah" becomes the active/selected sheet,
# zoomed to 49%, with E8 selected. Activate it last so it ends up as the
# workbook's active (tab-selected) sheet.
$ws1.Activate()
$win = $excel.ActiveWindow
$win.Zoom = 49
$ws1.Range("E8").Select()
